# CCC19 Derived Variables Spreadsheet — add derived variables for expected
# due dates of follow-up forms (30d / 90d / 180d), inserted (sorted by
# Variable #) just before the existing "X7" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above the current row 141 ("X7" / breast_biomarkers),
# shifting everything below down by three rows.
$ws.Rows("141:143").Insert()

# New row 141: X6a / 30d_due
$ws.Range("A141").Value = "X6a"
$ws.Range("B141").Value = "30d_due"
$ws.Range("C141").Value = "Other"
$ws.Range("D141").Value = "30d follow-up is due at the latest on this date"

# New row 142: X6b / 90d_due
$ws.Range("A142").Value = "X6b"
$ws.Range("B142").Value = "90d_due"
$ws.Range("C142").Value = "Other"
$ws.Range("D142").Value = "90d follow-up is due at the latest on this date"

# New row 143: X6c / 180d_due
$ws.Range("A143").Value = "X6c"
$ws.Range("B143").Value = "180d_due"
$ws.Range("C143").Value = "Other"
$ws.Range("D143").Value = "180d follow-up is due at the latest on this date"

# Grow Table1 to cover the three newly-inserted rows (was A1:E166, now A1:E169).
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:E169"))

# Restore view state: active cell / top-left cell scrolled to match.
$excel.ActiveWindow.ScrollRow = 129
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A143").Select()
